# Remove row 667 ("この広い宇宙に…" post) entirely; subsequent rows shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("667:667").Delete()
